# cursor position when moving between rows
#
# The task "when moving a row, move cursor to the row that is now in the
# position that the moved row previously held" (Id 32) is moved from the
# Active sheet to the Inactive sheet: its Status becomes "Done" and the
# "Done" (completed) date column is stamped with the same date it already
# carried ("4/11/2018"). All other rows shift accordingly.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Capture the row to move (Active row 2 = Id 32) before it shifts away.
$id = $active.Range("A2").Value2
$title = $active.Range("B2").Value2
$category = $active.Range("D2").Value2
$created = $active.Range("E2").Value2

# Remove it from the Active sheet; subsequent rows shift up.
$active.Rows.Item(2).Delete()

# Make room at the top of the Inactive sheet (below the header) and move
# the task there, now marked Done, with its completion date recorded.
$inactive.Rows.Item(2).Insert()

$inactive.Range("A2").Value = $id
$inactive.Range("B2").Value = "'" + $title
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "'" + $category
$inactive.Range("E2").Value = "'" + $created
$inactive.Range("F2").Value = "'4/11/2018"

$inactive.Range("A2:F2").Style = "Normal"
